$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "name" column (D) duplicating the "id" column (C) values for each field row.
$ws.Range("D1").Value2 = "name"
$ws.Range("D2").Value2 = $ws.Range("C2").Value2
$ws.Range("D3").Value2 = $ws.Range("C3").Value2
$ws.Range("D4").Value2 = $ws.Range("C4").Value2
$ws.Range("D5").Value2 = $ws.Range("C5").Value2
$ws.Range("D6").Value2 = $ws.Range("C6").Value2
$ws.Range("D7").Value2 = $ws.Range("C7").Value2

# Fill in min/max validation values for the Age row.
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 200

# Correct the example "items" values for the Profession (select) row (keep the quote-prefix cell style).
$ws.Range("K6").Value = "'A,b, C"

# Correct the example "items" values for the Gender (radio) row and remove the stray "-" in the max column.
$ws.Range("J7").ClearContents()
$ws.Range("K7").Value = "'Male, Female, Other"
